$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Uniform - Random")
$ws.Activate()

# Update the numeric values in columns C, D, E for rows 3-6, and C7
$ws.Range("C4").Value = 44699
$ws.Range("C5").Value = 44706
$ws.Range("C6").Value = 44698
$ws.Range("C7").Value = 44762

$ws.Range("D3").Value = 29452
$ws.Range("D4").Value = 29528
$ws.Range("D5").Value = 29658
$ws.Range("D6").Value = 29472

$ws.Range("E3").Value = 75800
$ws.Range("E4").Value = 75935
$ws.Range("E5").Value = 76197
$ws.Range("E6").Value = 76085

# Update the selected cell to match the saved view state
$ws.Range("C7").Select()
